$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.178715
$ws.Range("H2").Value = 0.536145
$ws.Range("I2").Value = 0.09904930989061336
$ws.Range("J2").Value = 0.09904930989061336
$ws.Range("M2").Value = 11.651608
$ws.Range("N2").Value = 34.954824
$ws.Range("O2").Value = 0.1892813629236475
$ws.Range("P2").Value = 0.1892813629236474
$ws.Range("Q2").Value = 2.08231712372
$ws.Range("R2").Value = 18.74085411348
$ws.Range("S2").Value = 0.01874818837274201
$ws.Range("T2").Value = 0.01874818837274201
$ws.Range("G3").Value = 0.178715
$ws.Range("H3").Value = 0.536145
$ws.Range("I3").Value = 0.09904930989061336
$ws.Range("J3").Value = 0.09904930989061336
$ws.Range("O3").Value = 0.4419371310876561
$ws.Range("P3").Value = 0.4419371310876561
$ws.Range("Q3").Value = 4.861827078256666
$ws.Range("R3").Value = 43.75644370431
$ws.Range("S3").Value = 0.04377356784926987
$ws.Range("T3").Value = 0.04377356784926987
$ws.Range("G4").Value = 0.178715
$ws.Range("H4").Value = 0.536145
$ws.Range("I4").Value = 0.09904930989061336
$ws.Range("J4").Value = 0.09904930989061336
$ws.Range("M4").Value = 8.657178999999999
$ws.Range("N4").Value = 25.971537
$ws.Range("O4").Value = 0.1406366091439035
$ws.Range("P4").Value = 0.1406366091439035
$ws.Range("Q4").Value = 1.547167744985
$ws.Range("R4").Value = 13.924509704865
$ws.Range("S4").Value = 0.01392995908105957
$ws.Range("T4").Value = 0.01392995908105957
$ws.Range("G5").Value = 0.178715
$ws.Range("H5").Value = 0.536145
$ws.Range("I5").Value = 0.09904930989061336
$ws.Range("J5").Value = 0.09904930989061336
$ws.Range("M5").Value = 5.488499666666667
$ws.Range("N5").Value = 16.465499
$ws.Range("O5").Value = 0.08916114387925267
$ws.Range("P5").Value = 0.08916114387925267
$ws.Range("Q5").Value = 0.9808772179283333
$ws.Range("R5").Value = 8.827894961355
$ws.Range("S5").Value = 0.008831349770297663
$ws.Range("T5").Value = 0.008831349770297663
$ws.Range("G6").Value = 0.178715
$ws.Range("H6").Value = 0.536145
$ws.Range("I6").Value = 0.09904930989061336
$ws.Range("J6").Value = 0.09904930989061336
$ws.Range("M6").Value = 4.091608333333333
$ws.Range("N6").Value = 12.274825
$ws.Range("O6").Value = 0.06646852536431769
$ws.Range("P6").Value = 0.06646852536431769
$ws.Range("Q6").Value = 0.7312317832916666
$ws.Range("R6").Value = 6.581086049625
$ws.Range("S6").Value = 0.006583661566782397
$ws.Range("T6").Value = 0.006583661566782397
$ws.Range("G7").Value = 0.178715
$ws.Range("H7").Value = 0.536145
$ws.Range("I7").Value = 0.09904930989061336
$ws.Range("J7").Value = 0.09904930989061336
$ws.Range("M7").Value = 4.463825666666667
$ws.Range("N7").Value = 13.391477
$ws.Range("O7").Value = 0.07251522760122259
$ws.Range("P7").Value = 0.07251522760122257
$ws.Range("Q7").Value = 0.7977526040183333
$ws.Range("R7").Value = 7.179773436165
$ws.Range("S7").Value = 0.007182583250461855
$ws.Range("T7").Value = 0.007182583250461855
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.450498333333333
$ws.Range("H8").Value = 4.351495
$ws.Range("I8").Value = 0.8039104659046613
$ws.Range("J8").Value = 0.8039104659046612
$ws.Range("M8").Value = 11.651608
$ws.Range("N8").Value = 34.954824
$ws.Range("O8").Value = 0.1892813629236475
$ws.Range("P8").Value = 0.1892813629236474
$ws.Range("Q8").Value = 16.90063798465333
$ws.Range("R8").Value = 152.10574186188
$ws.Range("S8").Value = 0.1521652686550187
$ws.Range("T8").Value = 0.1521652686550187
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.450498333333333
$ws.Range("H9").Value = 4.351495
$ws.Range("I9").Value = 0.8039104659046613
$ws.Range("J9").Value = 0.8039104659046612
$ws.Range("O9").Value = 0.4419371310876561
$ws.Range("P9").Value = 0.4419371310876561
$ws.Range("Q9").Value = 39.45987787240111
$ws.Range("R9").Value = 355.13890085161
$ws.Range("S9").Value = 0.355277884953247
$ws.Range("T9").Value = 0.3552778849532469
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.450498333333333
$ws.Range("H10").Value = 4.351495
$ws.Range("I10").Value = 0.8039104659046613
$ws.Range("J10").Value = 0.8039104659046612
$ws.Range("M10").Value = 8.657178999999999
$ws.Range("N10").Value = 25.971537
$ws.Range("O10").Value = 0.1406366091439035
$ws.Range("P10").Value = 0.1406366091439035
$ws.Range("Q10").Value = 12.55722371086833
$ws.Range("R10").Value = 113.015013397815
$ws.Range("S10").Value = 0.1130592419801272
$ws.Range("T10").Value = 0.1130592419801272
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.450498333333333
$ws.Range("H11").Value = 4.351495
$ws.Range("I11").Value = 0.8039104659046613
$ws.Range("J11").Value = 0.8039104659046612
$ws.Range("M11").Value = 5.488499666666667
$ws.Range("N11").Value = 16.465499
$ws.Range("O11").Value = 0.08916114387925267
$ws.Range("P11").Value = 0.08916114387925267
$ws.Range("Q11").Value = 7.961059619000555
$ws.Range("R11").Value = 71.649536571005
$ws.Range("S11").Value = 0.07167757671656255
$ws.Range("T11").Value = 0.07167757671656254
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.450498333333333
$ws.Range("H12").Value = 4.351495
$ws.Range("I12").Value = 0.8039104659046613
$ws.Range("J12").Value = 0.8039104659046612
$ws.Range("M12").Value = 4.091608333333333
$ws.Range("N12").Value = 12.274825
$ws.Range("O12").Value = 0.06646852536431769
$ws.Range("P12").Value = 0.06646852536431769
$ws.Range("Q12").Value = 5.934871068152777
$ws.Range("R12").Value = 53.413839613375
$ws.Range("S12").Value = 0.05343474319362443
$ws.Range("T12").Value = 0.05343474319362442
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.450498333333333
$ws.Range("H13").Value = 4.351495
$ws.Range("I13").Value = 0.8039104659046613
$ws.Range("J13").Value = 0.8039104659046612
$ws.Range("M13").Value = 4.463825666666667
$ws.Range("N13").Value = 13.391477
$ws.Range("O13").Value = 0.07251522760122259
$ws.Range("P13").Value = 0.07251522760122257
$ws.Range("Q13").Value = 6.474771689790555
$ws.Range("R13").Value = 58.272945208115
$ws.Range("S13").Value = 0.0582957504060814
$ws.Range("T13").Value = 0.05829575040608138
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.17509
$ws.Range("H14").Value = 0.52527
$ws.Range("I14").Value = 0.09704022420472538
$ws.Range("J14").Value = 0.09704022420472537
$ws.Range("M14").Value = 11.651608
$ws.Range("N14").Value = 34.954824
$ws.Range("O14").Value = 0.1892813629236475
$ws.Range("P14").Value = 0.1892813629236474
$ws.Range("Q14").Value = 2.04008004472
$ws.Range("R14").Value = 18.36072040248
$ws.Range("S14").Value = 0.01836790589588674
$ws.Range("T14").Value = 0.01836790589588674
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.17509
$ws.Range("H15").Value = 0.52527
$ws.Range("I15").Value = 0.09704022420472538
$ws.Range("J15").Value = 0.09704022420472537
$ws.Range("O15").Value = 0.4419371310876561
$ws.Range("P15").Value = 0.4419371310876561
$ws.Range("Q15").Value = 4.763211275673333
$ws.Range("R15").Value = 42.86890148106
$ws.Range("S15").Value = 0.04288567828513926
$ws.Range("T15").Value = 0.04288567828513926
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.17509
$ws.Range("H16").Value = 0.52527
$ws.Range("I16").Value = 0.09704022420472538
$ws.Range("J16").Value = 0.09704022420472537
$ws.Range("M16").Value = 8.657178999999999
$ws.Range("N16").Value = 25.971537
$ws.Range("O16").Value = 0.1406366091439035
$ws.Range("P16").Value = 0.1406366091439035
$ws.Range("Q16").Value = 1.51578547111
$ws.Range("R16").Value = 13.64206923999
$ws.Range("S16").Value = 0.01364740808271673
$ws.Range("T16").Value = 0.01364740808271673
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.17509
$ws.Range("H17").Value = 0.52527
$ws.Range("I17").Value = 0.09704022420472538
$ws.Range("J17").Value = 0.09704022420472537
$ws.Range("M17").Value = 5.488499666666667
$ws.Range("N17").Value = 16.465499
$ws.Range("O17").Value = 0.08916114387925267
$ws.Range("P17").Value = 0.08916114387925267
$ws.Range("Q17").Value = 0.9609814066366666
$ws.Range("R17").Value = 8.648832659730001
$ws.Range("S17").Value = 0.008652217392392457
$ws.Range("T17").Value = 0.008652217392392456
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 0.3333333333333333
$ws.Range("G18").Value = 0.17509
$ws.Range("H18").Value = 0.52527
$ws.Range("I18").Value = 0.09704022420472538
$ws.Range("J18").Value = 0.09704022420472537
$ws.Range("M18").Value = 4.091608333333333
$ws.Range("N18").Value = 12.274825
$ws.Range("O18").Value = 0.06646852536431769
$ws.Range("P18").Value = 0.06646852536431769
$ws.Range("Q18").Value = 0.7163997030833333
$ws.Range("R18").Value = 6.44759732775
$ws.Range("S18").Value = 0.006450120603910864
$ws.Range("T18").Value = 0.006450120603910863
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0.3333333333333333
$ws.Range("G19").Value = 0.17509
$ws.Range("H19").Value = 0.52527
$ws.Range("I19").Value = 0.09704022420472538
$ws.Range("J19").Value = 0.09704022420472537
$ws.Range("M19").Value = 4.463825666666667
$ws.Range("N19").Value = 13.391477
$ws.Range("O19").Value = 0.07251522760122259
$ws.Range("P19").Value = 0.07251522760122259
$ws.Range("Q19").Value = 0.7815712359766667
$ws.Range("R19").Value = 7.03414112379
$ws.Range("S19").Value = 0.00703689394467933
$ws.Range("T19").Value = 0.007036893944679328

Write-Output "Updated 234 cells in Bmp10-Acvr2a LR-pairs sheet with new TPM-based values."